# Update Sheets via scheduled runner - numeric profit/price recalculation
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 51099.617
$ws.Range("I100").Value = 59299.547
$ws.Range("J100").Value = 6000
$ws.Range("K100").Value = 59299.547
$ws.Range("L100").Value = 6000
$ws.Range("M100").Value = -58758.547
$ws.Range("N100").Value = -7082
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 2500
$ws.Range("K113").Value = 2500
$ws.Range("M113").Value = 754
$ws.Range("H132").Value = 5983801
$ws.Range("I132").Value = 6243927.5
$ws.Range("K132").Value = 18731782.5
$ws.Range("M132").Value = -18729252.5
$ws.Range("H141").Value = 5655.357
$ws.Range("I141").Value = 5606.8184
$ws.Range("K141").Value = 16820.4552
$ws.Range("M141").Value = -11640.4552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 7753.8
$ws.Range("I28").Value = 7753.8
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 7753.8
$ws.Range("L28").Value = 0
$ws.Range("M28").Value = -7561.8
$ws.Range("N28").ClearContents()
$ws.Range("H32").Value = 22400.81
$ws.Range("I32").Value = 22400.81
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 22400.81
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -22113.81
$ws.Range("N32").ClearContents()
$ws.Range("H61").Value = 4636.6
$ws.Range("I61").Value = 1487.9333
$ws.Range("J61").Value = 9359.6
$ws.Range("K61").Value = 1487.9333
$ws.Range("L61").Value = 9359.6
$ws.Range("M61").Value = -1275.9333
$ws.Range("N61").Value = -9783.6
$ws.Range("H99").Value = 7753.8
$ws.Range("I99").Value = 7753.8
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 7753.8
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -4758.8
$ws.Range("N99").ClearContents()
$ws.Range("H102").Value = 3743.4827
$ws.Range("I102").Value = 3291.1482
$ws.Range("K102").Value = 3291.1482
$ws.Range("M102").Value = -1669.1482
$ws.Range("H132").Value = 1516.6296
$ws.Range("I132").Value = 1244.1578
$ws.Range("K132").Value = 3732.4734
$ws.Range("M132").Value = -1202.4734
$ws.Range("H136").Value = 4636.6
$ws.Range("I136").Value = 1487.9333
$ws.Range("J136").Value = 9359.6
$ws.Range("K136").Value = 4463.7999
$ws.Range("L136").Value = 28078.8
$ws.Range("M136").Value = -1913.7999
$ws.Range("N136").Value = -33178.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 95568.375
$ws.Range("J132").Value = 95568.375
$ws.Range("L132").Value = 95568.375
$ws.Range("N132").Value = -105688.375
$ws.Range("H134").Value = 2384.932
$ws.Range("J134").Value = 4197
$ws.Range("L134").Value = 12591
$ws.Range("N134").Value = -17661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5265255.5
$ws.Range("I31").Value = 6668011.5
$ws.Range("K31").Value = 6668011.5
$ws.Range("M31").Value = -6667716.5
$ws.Range("H34").Value = 5265255.5
$ws.Range("I34").Value = 6668011.5
$ws.Range("K34").Value = 6668011.5
$ws.Range("M34").Value = -6667809.5
$ws.Range("H99").Value = 8419.666999999999
$ws.Range("I99").Value = 6858.6665
$ws.Range("K99").Value = 6858.6665
$ws.Range("M99").Value = -5360.6665
$ws.Range("H105").Value = 1888.25
$ws.Range("I105").Value = 1132.5
$ws.Range("K105").Value = 1132.5
$ws.Range("M105").Value = 614.5
$ws.Range("H126").Value = 8419.666999999999
$ws.Range("I126").Value = 6858.6665
$ws.Range("K126").Value = 20575.9995
$ws.Range("M126").Value = -18105.9995
$ws.Range("H134").Value = 2719.8
$ws.Range("I134").Value = 2599.7693
$ws.Range("K134").Value = 7799.3079
$ws.Range("M134").Value = -5264.3079

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 39735530
$ws.Range("I4").Value = 49337430
$ws.Range("K4").Value = 148012290
$ws.Range("M4").Value = -148012178
$ws.Range("H11").Value = 85198.836
$ws.Range("I11").Value = 126098.625
$ws.Range("J11").Value = 3399.25
$ws.Range("K11").Value = 378295.875
$ws.Range("L11").Value = 10197.75
$ws.Range("M11").Value = -378155.875
$ws.Range("N11").Value = -10477.75
$ws.Range("H17").Value = 165.09091
$ws.Range("I17").Value = 129.33333
$ws.Range("K17").Value = 387.99999
$ws.Range("M17").Value = -218.99999
$ws.Range("H39").Value = 3675
$ws.Range("J39").Value = 3675
$ws.Range("L39").Value = 11025
$ws.Range("N39").Value = -11613
$ws.Range("H44").Value = 2947.3
$ws.Range("I44").Value = 1886.5
$ws.Range("J44").Value = 3212.5
$ws.Range("K44").Value = 5659.5
$ws.Range("L44").Value = 9637.5
$ws.Range("M44").Value = -5261.5
$ws.Range("N44").Value = -10433.5
$ws.Range("H55").Value = 3211.875
$ws.Range("J55").Value = 4666.6665
$ws.Range("L55").Value = 13999.9995
$ws.Range("N55").Value = -14353.9995
$ws.Range("H131").Value = 1841.6
$ws.Range("I131").Value = 1753.1666
$ws.Range("K131").Value = 5259.4998
$ws.Range("M131").Value = -219.4997999999996
$ws.Range("H132").Value = 3699.5
$ws.Range("I132").Value = 3699.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 33295.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -30765.5
$ws.Range("N132").ClearContents()
$ws.Range("H139").Value = 4847.1055
$ws.Range("I139").Value = 4977.5
$ws.Range("K139").Value = 14932.5
$ws.Range("M139").Value = -9792.5
$ws.Range("H140").Value = 5360.4443
$ws.Range("I140").Value = 4781.125
$ws.Range("K140").Value = 14343.375
$ws.Range("M140").Value = -9163.375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1413.8334
$ws.Range("I113").Value = 662.8889
$ws.Range("J113").Value = 3666.6667
$ws.Range("K113").Value = 662.8889
$ws.Range("L113").Value = 3666.6667
$ws.Range("M113").Value = 1507.1111
$ws.Range("N113").Value = -8006.6667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 916.875
$ws.Range("I22").Value = 583.375
$ws.Range("K22").Value = 583.375
$ws.Range("M22").Value = -288.375
$ws.Range("H27").Value = 916.875
$ws.Range("I27").Value = 583.375
$ws.Range("K27").Value = 583.375
$ws.Range("M27").Value = -476.375
$ws.Range("H46").Value = 5798.5557
$ws.Range("I46").Value = 1749.8334
$ws.Range("J46").Value = 7822.9165
$ws.Range("K46").Value = 1749.8334
$ws.Range("L46").Value = 7822.9165
$ws.Range("M46").Value = -1561.8334
$ws.Range("N46").Value = -8198.916499999999
$ws.Range("H68").Value = 4266.696
$ws.Range("I68").Value = 3426.8462
$ws.Range("K68").Value = 3426.8462
$ws.Range("M68").Value = -2677.8462
$ws.Range("H71").Value = 4266.696
$ws.Range("I71").Value = 3426.8462
$ws.Range("K71").Value = 17134.231
$ws.Range("M71").Value = -13390.231
$ws.Range("H93").Value = 3113.8635
$ws.Range("I93").Value = 2659.2354
$ws.Range("K93").Value = 2659.2354
$ws.Range("M93").Value = -1411.2354
$ws.Range("H100").Value = 2321.875
$ws.Range("I100").Value = 2321.875
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 2321.875
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -1780.875
$ws.Range("N100").ClearContents()
$ws.Range("H132").Value = 6176.15
$ws.Range("I132").Value = 6048.6
$ws.Range("K132").Value = 18145.8
$ws.Range("M132").Value = -15615.8
$ws.Range("H136").Value = 5816.857
$ws.Range("I136").Value = 5828
$ws.Range("K136").Value = 17484
$ws.Range("M136").Value = -14934

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 890
$ws.Range("I81").Value = 890
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1780
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -719
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 890
$ws.Range("I84").Value = 890
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 8900
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -3596
$ws.Range("N84").ClearContents()
$ws.Range("H107").Value = 2292.25
$ws.Range("I107").Value = 1934.5
$ws.Range("K107").Value = 5803.5
$ws.Range("M107").Value = -3883.5
$ws.Range("H126").Value = 3444.342
$ws.Range("I126").Value = 3142.5186
$ws.Range("K126").Value = 9427.5558
$ws.Range("M126").Value = -6957.5558
$ws.Range("H136").Value = 13014.55
$ws.Range("I136").Value = 14851.327
$ws.Range("K136").Value = 44553.981
$ws.Range("M136").Value = -42003.981
